# Typo fix in URL
#
# The "Working Papers" sheet has a hyperlinked cell (D4) whose displayed
# text is the GE-IMR working-paper URL. It contained a typo
# ("ggkilleen33" instead of "gkilleen33"); fix the cell text so it reads
# the correct URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Working Papers")

$cell = $ws.Range("D4")
$cell.Value = "https://gkilleen33.github.io/papers/working/GE-IMR.pdf"
